# Tutorial 6 solution update: switch date strings from DD/MM/YYYY to
# DD-MM-YYYY and update the derived attendance counters for the rows whose
# totals changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Date text, D, E, F, G, H)
$rows = @(
    @{ Row = 3;  Date = "28-07-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 },
    @{ Row = 4;  Date = "01-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 5;  Date = "04-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 6;  Date = "08-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 7;  Date = "11-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 8;  Date = "15-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 9;  Date = "18-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 10; Date = "22-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 11; Date = "25-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 12; Date = "29-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 13; Date = "01-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 14; Date = "05-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 15; Date = "08-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 16; Date = "12-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 17; Date = "15-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 18; Date = "19-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 19; Date = "22-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 20; Date = "26-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 21; Date = "29-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r.Row, 1)
    # Some of the new DD-MM-YYYY strings are ambiguous (day <= 12) and would
    # otherwise be auto-parsed into a date serial by the smart-entry logic.
    # Force plain text, then strip the temporary "Text" number format so the
    # cell's style matches an untouched cell again.
    $cell.NumberFormat = "@"
    $cell.Value = $r.Date
    $cell.ClearFormats()

    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}
